$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.507.99'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.14%  '
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.160.00'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -2.41%  '
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.14'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.71%  '
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.604'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -3.38%  '
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '70.96'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.87%  '
$ws.Range("E7").Style = "Normal"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.573'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -4.55%  '
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.79'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -5.19%  '
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0901'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -5.04%  '
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.20'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -4.69%  '
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.100'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -3.44%  '
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.65'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -4.56%  '
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.479.84'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.57%  '
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.30'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.78%  '
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.139.63'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -3.18%  '
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.783'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -5.93%  '
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '41.321.90'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.17%  '
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -4.31%  '
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '69.31'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -4.19%  '
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.74'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -6.47%  '
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.69'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -11.55%  '
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '225.51'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.58%  '
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.98'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.95%  '
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.39%  '
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.63'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -6.95%  '
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.29'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -9.22%  '
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.17'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -3.89%  '
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.17'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.51%  '
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '171.32'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +2.66%  '
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '19.74'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.24%  '
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '32.15'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +6.93%  '
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0760'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -4.06%  '
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.08'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -9.46%  '
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -3.56%  '
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.28'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.84%  '
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.104'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.20%  '
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0292'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.45%  '
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '12.11'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -12.38%  '
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.06'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.43%  '
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.29'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -5.61%  '
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '58.51'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -9.03%  '
$ws.Range("E43").Style = "Normal"

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -3.98%  '
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.25'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -4.46%  '
$ws.Range("E45").Style = "Normal"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -3.89%  '
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '97.57'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -5.52%  '
$ws.Range("E47").Style = "Normal"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.18%  '
$ws.Range("E48").Style = "Normal"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -4.85%  '
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.16'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -7.62%  '
$ws.Range("E50").Style = "Normal"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -3.06%  '
$ws.Range("E51").Style = "Normal"
